# h5getInfo.xlsx update:
#  - Re-format the "Shape:" value on the h5getInfo1 sheet from "(721, 841)"
#    to "{ 721, 841 }" (escaping ',' in compound-type field names).
#  - Make "h5getInfo1" the active sheet/tab (was "Globals").

$wb = $excel.ActiveWorkbook

$h5getInfo1 = $wb.Worksheets.Item("h5getInfo1")

# Update the Shape: value (row 8, column B) with the new formatting.
$h5getInfo1.Range("B8").Value = "{ 721, 841 }"

# Switch the active sheet from Globals to h5getInfo1, with A2 selected
# (matching the selection previously shown on the Globals sheet).
$h5getInfo1.Activate()
$h5getInfo1.Range("A2").Select()
